$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Init")

# Update the Lower Right Cell value for the Variables List Indices row (row 15)
# from "E238" to "E239" (extends the variables listing block by one row).
$ws.Range("D15").Value = "E239"

# Reflect the new active selection on the sheet.
$ws.Activate()
$ws.Range("D15").Select()
